$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 / A20: the recorded timestamp shifted slightly on re-fetch.
$ws.Cells.Item(20, 1).Value = 44333.78536134028

# New data row retrieved on 2021-05-18 18:51:13 UTC
$row = 21
$ws.Cells.Item($row, 1).Value  = 44334.78557194959
$ws.Cells.Item($row, 2).Value  = 73902
$ws.Cells.Item($row, 3).Value  = 62176
$ws.Cells.Item($row, 4).Value  = 3349
$ws.Cells.Item($row, 5).Value  = 2076
$ws.Cells.Item($row, 6).Value  = 1467
$ws.Cells.Item($row, 7).Value  = 19215
$ws.Cells.Item($row, 8).Value  = 1382
$ws.Cells.Item($row, 9).Value  = 847
$ws.Cells.Item($row, 10).Value = 214

# Keep the date-formatted style consistent with the rest of column A.
$ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item($row - 1, 1).NumberFormat
